$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing activity log entry text in G40 (shared string used there)
$ws.Range("G40").Value = "Re-compiled VHDL code"

# Fill E40 (end time for row 40)
$ws.Range("E40").Value = 0.020833333333333332

# Fill row 41 entries
$ws.Range("B41").Value = 6977
$ws.Range("C41").Value = 43926
$ws.Range("D41").Value = 0.020833333333333332
$ws.Range("E41").Value = 0.036111111111111115
$ws.Range("G41").Value = "Updated Functional Waveforms for LogicUnit.vhd"

# Update view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("G43").Select()
